$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("루밍")

# Correct the passenger/passport name in the comparison list:
# "YU/HYANGSUK" -> "YU/HYANGSUKE"
$ws.Range("B2").Value = "YU/HYANGSUKE"
